# Updates odds-database rows: for a handful of fixture pairs that share the
# same match date, the two rows had been entered swapped relative to the
# canonical (id-ordered) source feed. This re-synchronizes each pair by
# swapping the full data payload (columns B..AC) between the two rows,
# while leaving column A (the sequential row id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (1-based worksheet rows) whose B:AC contents must be swapped.
$rowPairs = @(
    @(194, 195),
    @(214, 215),
    @(221, 222),
    @(227, 228),
    @(248, 249),
    @(258, 259),
    @(264, 265)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # NOTE: use ${r1}/${r2} (braced) rather than bare $r1:AC$r1 -- a bare
    # "$r1:AC" inside a double-quoted string is parsed as a scope/drive
    # qualifier (like $env:VAR) and silently truncates the address.
    $range1 = $ws.Range("B${r1}:AC${r1}")
    $range2 = $ws.Range("B${r2}:AC${r2}")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
